$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:F2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("J2:K2").NumberFormat = "@"

$ws.Range("E2").Value = "128.20000000"
$ws.Range("F2").Value = "44560.50574000"
$ws.Range("H2").Value = "5709439.42501910"
$ws.Range("I2").Value = 22611
$ws.Range("J2").Value = "21721.68202000"
$ws.Range("K2").Value = "2782972.72016410"
$ws.Range("M2").Value = 128.2
$ws.Range("R2").Value = 128.2
$ws.Range("S2").Value = 128.2
$ws.Range("T2").Value = 128.2
$ws.Range("U2").Value = 128.2
